$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1663, 3639, 3007, 1, 936, 13574, 4621, 7416, 5447, 2494, 3009, 2576, 6357, 428, 8274, 6140, 863, 18900, 3388, 14055, 10422, 1950, 4284, 17941, 2136, 8133, 3227, 7579, 5307, 4308, 1523, 11192, 2472, 3848, 2503, 6978, 11359, 9826, 2784, 1246, 3910, 5813, 5179, 1963, 3053, 7602, 5552, 4148, 775, 1947, 0, 2442, 8897, 2620, 7119, 4181, 1920, 17532, 13817, 3387, 1279, 11738, 1460, 11844, 5571, 4896, 4156, 4014, 814, 8579, 2699, 3319, 13014, 13495, 1344, 1564, 3362, 3763, 2723, 3477, 3973, 4559, 1056, 16909, 5769, 1750, 6915, 3767, 10940, 11179, 1047, 2960, 3879, 6176, 3915, 7549, 0)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 4).Value = $v
    $row++
}

